$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Cells.Item(32, 8).Value = 2976.5217   # H32: 4284.5 -> 2976.5217
$ws.Cells.Item(32, 9).Value = 1280   # I32: 2880.2 -> 1280
$ws.Cells.Item(32, 10).Value = 3881.3333   # J32: 5064.6665 -> 3881.3333
$ws.Cells.Item(32, 11).Value = 1280   # K32: 2880.2 -> 1280
$ws.Cells.Item(32, 12).Value = 3881.3333   # L32: 5064.6665 -> 3881.3333
$ws.Cells.Item(32, 13).Value = -954   # M32: -2554.2 -> -954
$ws.Cells.Item(32, 14).Value = -4533.3333   # N32: -5716.6665 -> -4533.3333
# Row 64
$ws.Cells.Item(64, 8).Value = 3185.1538   # H64: 3011.5667 -> 3185.1538
$ws.Cells.Item(64, 9).Value = 2999.6   # I64: 2883 -> 2999.6
$ws.Cells.Item(64, 10).Value = 3301.125   # J64: 3043.7083 -> 3301.125
$ws.Cells.Item(64, 11).Value = 2999.6   # K64: 2883 -> 2999.6
$ws.Cells.Item(64, 12).Value = 3301.125   # L64: 3043.7083 -> 3301.125
$ws.Cells.Item(64, 13).Value = -2751.6   # M64: -2635 -> -2751.6
$ws.Cells.Item(64, 14).Value = -3797.125   # N64: -3539.7083 -> -3797.125
# Row 67
$ws.Cells.Item(67, 8).Value = 3185.1538   # H67: 3011.5667 -> 3185.1538
$ws.Cells.Item(67, 9).Value = 2999.6   # I67: 2883 -> 2999.6
$ws.Cells.Item(67, 10).Value = 3301.125   # J67: 3043.7083 -> 3301.125
$ws.Cells.Item(67, 11).Value = 2999.6   # K67: 2883 -> 2999.6
$ws.Cells.Item(67, 12).Value = 3301.125   # L67: 3043.7083 -> 3301.125
$ws.Cells.Item(67, 13).Value = -2141.6   # M67: -2025 -> -2141.6
$ws.Cells.Item(67, 14).Value = -5017.125   # N67: -4759.7083 -> -5017.125
# Row 82
$ws.Cells.Item(82, 8).Value = 747.1429000000001   # H82: 1126.25 -> 747.1429000000001
$ws.Cells.Item(82, 9).Value = 747.1429000000001   # I82: 744.2857 -> 747.1429000000001
$ws.Cells.Item(82, 10).Value = 0   # J82: 3800 -> 0
$ws.Cells.Item(82, 11).Value = 2241.4287   # K82: 2232.8571 -> 2241.4287
$ws.Cells.Item(82, 12).Value = 0   # L82: 11400 -> 0
$ws.Cells.Item(82, 13).Value = -1835.4287   # M82: -1826.8571 -> -1835.4287
$ws.Cells.Item(82, 14).Value = $null   # N82: -12212 -> None
# Row 85
$ws.Cells.Item(85, 8).Value = 747.1429000000001   # H85: 1126.25 -> 747.1429000000001
$ws.Cells.Item(85, 9).Value = 747.1429000000001   # I85: 744.2857 -> 747.1429000000001
$ws.Cells.Item(85, 10).Value = 0   # J85: 3800 -> 0
$ws.Cells.Item(85, 11).Value = 2241.4287   # K85: 2232.8571 -> 2241.4287
$ws.Cells.Item(85, 12).Value = 0   # L85: 11400 -> 0
$ws.Cells.Item(85, 13).Value = -837.4287000000004   # M85: -828.8571000000002 -> -837.4287000000004
$ws.Cells.Item(85, 14).Value = $null   # N85: -14208 -> None
# Row 113
$ws.Cells.Item(113, 8).Value = 2467.2812   # H113: 2909.7896 -> 2467.2812
$ws.Cells.Item(113, 9).Value = 1656.1111   # I113: 2651.25 -> 1656.1111
$ws.Cells.Item(113, 10).Value = 2784.6956   # J113: 2978.7334 -> 2784.6956
$ws.Cells.Item(113, 11).Value = 1656.1111   # K113: 2651.25 -> 1656.1111
$ws.Cells.Item(113, 12).Value = 2784.6956   # L113: 2978.7334 -> 2784.6956
$ws.Cells.Item(113, 13).Value = 1597.8889   # M113: 602.75 -> 1597.8889
$ws.Cells.Item(113, 14).Value = -9292.695599999999   # N113: -9486.733400000001 -> -9292.695599999999
# Row 116
$ws.Cells.Item(116, 8).Value = 6747989.5   # H116: 6441256 -> 6747989.5
$ws.Cells.Item(116, 9).Value = 7266906.5   # I116: 7085204 -> 7266906.5
$ws.Cells.Item(116, 10).Value = 2066.3333   # J116: 1776.25 -> 2066.3333
$ws.Cells.Item(116, 11).Value = 7266906.5   # K116: 7085204 -> 7266906.5
$ws.Cells.Item(116, 12).Value = 2066.3333   # L116: 1776.25 -> 2066.3333
$ws.Cells.Item(116, 13).Value = -7263464.5   # M116: -7081762 -> -7263464.5
$ws.Cells.Item(116, 14).Value = -8950.3333   # N116: -8660.25 -> -8950.3333
# Row 132
$ws.Cells.Item(132, 8).Value = 5110.6206   # H132: 4558.387 -> 5110.6206
$ws.Cells.Item(132, 9).Value = 5369.12   # I132: 4715.926 -> 5369.12
$ws.Cells.Item(132, 11).Value = 16107.36   # K132: 14147.778 -> 16107.36
$ws.Cells.Item(132, 13).Value = -13577.36   # M132: -11617.778 -> -13577.36

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Cells.Item(74, 8).Value = 12164567   # H74: 13237892 -> 12164567
$ws.Cells.Item(74, 9).Value = 16668809   # I74: 19567682 -> 16668809
$ws.Cells.Item(74, 10).Value = 3112.7   # J74: 2876.9092 -> 3112.7
$ws.Cells.Item(74, 11).Value = 16668809   # K74: 19567682 -> 16668809
$ws.Cells.Item(74, 12).Value = 3112.7   # L74: 2876.9092 -> 3112.7
$ws.Cells.Item(74, 13).Value = -16667935   # M74: -19566808 -> -16667935
$ws.Cells.Item(74, 14).Value = -4860.7   # N74: -4624.9092 -> -4860.7
# Row 77
$ws.Cells.Item(77, 8).Value = 12164567   # H77: 13237892 -> 12164567
$ws.Cells.Item(77, 9).Value = 16668809   # I77: 19567682 -> 16668809
$ws.Cells.Item(77, 10).Value = 3112.7   # J77: 2876.9092 -> 3112.7
$ws.Cells.Item(77, 11).Value = 83344045   # K77: 97838410 -> 83344045
$ws.Cells.Item(77, 12).Value = 15563.5   # L77: 14384.546 -> 15563.5
$ws.Cells.Item(77, 13).Value = -83339677   # M77: -97834042 -> -83339677
$ws.Cells.Item(77, 14).Value = -24299.5   # N77: -23120.546 -> -24299.5
# Row 97
$ws.Cells.Item(97, 8).Value = 3240   # H97: 15154867 -> 3240
$ws.Cells.Item(97, 9).Value = 3442.6365   # I97: 16670253 -> 3442.6365
$ws.Cells.Item(97, 11).Value = 3442.6365   # K97: 16670253 -> 3442.6365
$ws.Cells.Item(97, 13).Value = -2946.6365   # M97: -16669757 -> -2946.6365
# Row 102
$ws.Cells.Item(102, 8).Value = 2887.5   # H102: 90910900 -> 2887.5
$ws.Cells.Item(102, 9).Value = 2775   # I102: 142858270 -> 2775
$ws.Cells.Item(102, 11).Value = 2775   # K102: 142858270 -> 2775
$ws.Cells.Item(102, 13).Value = -1153   # M102: -142856648 -> -1153

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Cells.Item(86, 8).Value = 25002430   # H86: 22224838 -> 25002430
$ws.Cells.Item(86, 9).Value = 40001660   # I86: 40002240 -> 40001660
$ws.Cells.Item(86, 10).Value = 3717   # J86: 3087.75 -> 3717
$ws.Cells.Item(86, 11).Value = 40001660   # K86: 40002240 -> 40001660
$ws.Cells.Item(86, 12).Value = 3717   # L86: 3087.75 -> 3717
$ws.Cells.Item(86, 13).Value = -40000537   # M86: -40001117 -> -40000537
$ws.Cells.Item(86, 14).Value = -5963   # N86: -5333.75 -> -5963
# Row 89
$ws.Cells.Item(89, 8).Value = 25002430   # H89: 22224838 -> 25002430
$ws.Cells.Item(89, 9).Value = 40001660   # I89: 40002240 -> 40001660
$ws.Cells.Item(89, 10).Value = 3717   # J89: 3087.75 -> 3717
$ws.Cells.Item(89, 11).Value = 200008300   # K89: 200011200 -> 200008300
$ws.Cells.Item(89, 12).Value = 18585   # L89: 15438.75 -> 18585
$ws.Cells.Item(89, 13).Value = -200002684   # M89: -200005584 -> -200002684
$ws.Cells.Item(89, 14).Value = -29817   # N89: -26670.75 -> -29817
# Row 94
$ws.Cells.Item(94, 8).Value = 6121.6665   # H94: 6288.2856 -> 6121.6665
$ws.Cells.Item(94, 10).Value = 29031.715   # J94: 33822 -> 29031.715
$ws.Cells.Item(94, 12).Value = 29031.715   # L94: 33822 -> 29031.715
$ws.Cells.Item(94, 14).Value = -29933.715   # N94: -34724 -> -29933.715
# Row 99
$ws.Cells.Item(99, 8).Value = 1428   # H99: 1426 -> 1428
$ws.Cells.Item(99, 9).Value = 1299.9166   # I99: 1308.1818 -> 1299.9166
$ws.Cells.Item(99, 10).Value = 1647.5714   # J99: 1750 -> 1647.5714
$ws.Cells.Item(99, 11).Value = 1299.9166   # K99: 1308.1818 -> 1299.9166
$ws.Cells.Item(99, 12).Value = 1647.5714   # L99: 1750 -> 1647.5714
$ws.Cells.Item(99, 13).Value = 198.0834   # M99: 189.8181999999999 -> 198.0834
$ws.Cells.Item(99, 14).Value = -4643.5714   # N99: -4746 -> -4643.5714
# Row 105
$ws.Cells.Item(105, 8).Value = 3021.875   # H105: 2152.353 -> 3021.875
$ws.Cells.Item(105, 9).Value = 2645.8333   # I105: 1761.7273 -> 2645.8333
$ws.Cells.Item(105, 10).Value = 4150   # J105: 2868.5 -> 4150
$ws.Cells.Item(105, 11).Value = 2645.8333   # K105: 1761.7273 -> 2645.8333
$ws.Cells.Item(105, 12).Value = 4150   # L105: 2868.5 -> 4150
$ws.Cells.Item(105, 13).Value = -898.8332999999998   # M105: -14.72730000000001 -> -898.8332999999998
$ws.Cells.Item(105, 14).Value = -7644   # N105: -6362.5 -> -7644

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Cells.Item(62, 8).Value = 71431400   # H62: 2999.5 -> 71431400
$ws.Cells.Item(62, 9).Value = 3017.7273   # I62: 2999.5 -> 3017.7273
$ws.Cells.Item(62, 10).Value = 333335460   # J62: 0 -> 333335460
$ws.Cells.Item(62, 11).Value = 3017.7273   # K62: 2999.5 -> 3017.7273
$ws.Cells.Item(62, 12).Value = 333335460   # L62: 0 -> 333335460
$ws.Cells.Item(62, 13).Value = -2393.7273   # M62: -2375.5 -> -2393.7273
$ws.Cells.Item(62, 14).Value = -333336708   # N62: None -> -333336708
# Row 65
$ws.Cells.Item(65, 8).Value = 71431400   # H65: 2999.5 -> 71431400
$ws.Cells.Item(65, 9).Value = 3017.7273   # I65: 2999.5 -> 3017.7273
$ws.Cells.Item(65, 10).Value = 333335460   # J65: 0 -> 333335460
$ws.Cells.Item(65, 11).Value = 15088.6365   # K65: 14997.5 -> 15088.6365
$ws.Cells.Item(65, 12).Value = 1666677300   # L65: 0 -> 1666677300
$ws.Cells.Item(65, 13).Value = -11968.6365   # M65: -11877.5 -> -11968.6365
$ws.Cells.Item(65, 14).Value = -1666683540   # N65: None -> -1666683540
# Row 99
$ws.Cells.Item(99, 8).Value = 1883718   # H99: 1789544.1 -> 1883718
$ws.Cells.Item(99, 9).Value = 2980560.8   # I99: 3576453 -> 2980560.8
$ws.Cells.Item(99, 10).Value = 3416.1428   # J99: 2635.4 -> 3416.1428
$ws.Cells.Item(99, 11).Value = 2980560.8   # K99: 3576453 -> 2980560.8
$ws.Cells.Item(99, 12).Value = 3416.1428   # L99: 2635.4 -> 3416.1428
$ws.Cells.Item(99, 13).Value = -2979062.8   # M99: -3574955 -> -2979062.8
$ws.Cells.Item(99, 14).Value = -6412.1428   # N99: -5631.4 -> -6412.1428
# Row 126
$ws.Cells.Item(126, 8).Value = 1883718   # H126: 1789544.1 -> 1883718
$ws.Cells.Item(126, 9).Value = 2980560.8   # I126: 3576453 -> 2980560.8
$ws.Cells.Item(126, 10).Value = 3416.1428   # J126: 2635.4 -> 3416.1428
$ws.Cells.Item(126, 11).Value = 8941682.399999999   # K126: 10729359 -> 8941682.399999999
$ws.Cells.Item(126, 12).Value = 10248.4284   # L126: 7906.200000000001 -> 10248.4284
$ws.Cells.Item(126, 13).Value = -8939212.399999999   # M126: -10726889 -> -8939212.399999999
$ws.Cells.Item(126, 14).Value = -15188.4284   # N126: -12846.2 -> -15188.4284

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Cells.Item(131, 8).Value = 850.125   # H131: 5377169.5 -> 850.125
$ws.Cells.Item(131, 9).Value = 428.18182   # I131: 397.85715 -> 428.18182
$ws.Cells.Item(131, 10).Value = 975.56757   # J131: 6945394.5 -> 975.56757
$ws.Cells.Item(131, 11).Value = 1284.54546   # K131: 1193.57145 -> 1284.54546
$ws.Cells.Item(131, 12).Value = 2926.70271   # L131: 20836183.5 -> 2926.70271
$ws.Cells.Item(131, 13).Value = 3755.45454   # M131: 3846.42855 -> 3755.45454
$ws.Cells.Item(131, 14).Value = -13006.70271   # N131: -20846263.5 -> -13006.70271

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Cells.Item(80, 8).Value = 3487.4443   # H80: 3624.625 -> 3487.4443
$ws.Cells.Item(80, 9).Value = 4128.846   # I80: 4455 -> 4128.846
$ws.Cells.Item(80, 10).Value = 2891.8572   # J80: 2922 -> 2891.8572
$ws.Cells.Item(80, 11).Value = 4128.846   # K80: 4455 -> 4128.846
$ws.Cells.Item(80, 12).Value = 2891.8572   # L80: 2922 -> 2891.8572
$ws.Cells.Item(80, 13).Value = -3130.846   # M80: -3457 -> -3130.846
$ws.Cells.Item(80, 14).Value = -4887.8572   # N80: -4918 -> -4887.8572
# Row 83
$ws.Cells.Item(83, 8).Value = 3487.4443   # H83: 3624.625 -> 3487.4443
$ws.Cells.Item(83, 9).Value = 4128.846   # I83: 4455 -> 4128.846
$ws.Cells.Item(83, 10).Value = 2891.8572   # J83: 2922 -> 2891.8572
$ws.Cells.Item(83, 11).Value = 20644.23   # K83: 22275 -> 20644.23
$ws.Cells.Item(83, 12).Value = 14459.286   # L83: 14610 -> 14459.286
$ws.Cells.Item(83, 13).Value = -15652.23   # M83: -17283 -> -15652.23
$ws.Cells.Item(83, 14).Value = -24443.286   # N83: -24594 -> -24443.286
# Row 97
$ws.Cells.Item(97, 8).Value = 1206.909   # H97: 1235.8096 -> 1206.909
$ws.Cells.Item(97, 9).Value = 1237.2222   # I97: 1274.7059 -> 1237.2222
$ws.Cells.Item(97, 11).Value = 1237.2222   # K97: 1274.7059 -> 1237.2222
$ws.Cells.Item(97, 13).Value = -741.2221999999999   # M97: -778.7058999999999 -> -741.2221999999999

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Cells.Item(82, 8).Value = 3047.762   # H82: 3977.1538 -> 3047.762
$ws.Cells.Item(82, 9).Value = 2081.25   # I82: 2688.889 -> 2081.25
$ws.Cells.Item(82, 10).Value = 6140.6   # J82: 6875.75 -> 6140.6
$ws.Cells.Item(82, 11).Value = 2081.25   # K82: 2688.889 -> 2081.25
$ws.Cells.Item(82, 12).Value = 6140.6   # L82: 6875.75 -> 6140.6
$ws.Cells.Item(82, 13).Value = -1720.25   # M82: -2327.889 -> -1720.25
$ws.Cells.Item(82, 14).Value = -6862.6   # N82: -7597.75 -> -6862.6
# Row 85
$ws.Cells.Item(85, 8).Value = 3047.762   # H85: 3977.1538 -> 3047.762
$ws.Cells.Item(85, 9).Value = 2081.25   # I85: 2688.889 -> 2081.25
$ws.Cells.Item(85, 10).Value = 6140.6   # J85: 6875.75 -> 6140.6
$ws.Cells.Item(85, 11).Value = 2081.25   # K85: 2688.889 -> 2081.25
$ws.Cells.Item(85, 12).Value = 6140.6   # L85: 6875.75 -> 6140.6
$ws.Cells.Item(85, 13).Value = -833.25   # M85: -1440.889 -> -833.25
$ws.Cells.Item(85, 14).Value = -8636.6   # N85: -9371.75 -> -8636.6
# Row 93
$ws.Cells.Item(93, 8).Value = 1280.3513   # H93: 27779816 -> 1280.3513
$ws.Cells.Item(93, 9).Value = 1183   # I93: 29413628 -> 1183
$ws.Cells.Item(93, 10).Value = 1783.3334   # J93: 5000 -> 1783.3334
$ws.Cells.Item(93, 11).Value = 1183   # K93: 29413628 -> 1183
$ws.Cells.Item(93, 12).Value = 1783.3334   # L93: 5000 -> 1783.3334
$ws.Cells.Item(93, 13).Value = 65   # M93: -29412380 -> 65
$ws.Cells.Item(93, 14).Value = -4279.3334   # N93: -7496 -> -4279.3334
# Row 100
$ws.Cells.Item(100, 8).Value = 1560.3   # H100: 1545.8182 -> 1560.3
$ws.Cells.Item(100, 9).Value = 1199.8   # I100: 1228.5714 -> 1199.8
$ws.Cells.Item(100, 10).Value = 1920.8   # J100: 2101 -> 1920.8
$ws.Cells.Item(100, 11).Value = 1199.8   # K100: 1228.5714 -> 1199.8
$ws.Cells.Item(100, 12).Value = 1920.8   # L100: 2101 -> 1920.8
$ws.Cells.Item(100, 13).Value = -658.8   # M100: -687.5714 -> -658.8
$ws.Cells.Item(100, 14).Value = -3002.8   # N100: -3183 -> -3002.8

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 4370.222   # H62: 111114980 -> 4370.222
$ws.Cells.Item(62, 9).Value = 3574.75   # I62: 0 -> 3574.75
$ws.Cells.Item(62, 10).Value = 5006.6   # J62: 111114980 -> 5006.6
$ws.Cells.Item(62, 11).Value = 3574.75   # K62: 0 -> 3574.75
$ws.Cells.Item(62, 12).Value = 5006.6   # L62: 111114980 -> 5006.6
$ws.Cells.Item(62, 13).Value = -2950.75   # M62: None -> -2950.75
$ws.Cells.Item(62, 14).Value = -6254.6   # N62: -111116228 -> -6254.6
# Row 65
$ws.Cells.Item(65, 8).Value = 4370.222   # H65: 111114980 -> 4370.222
$ws.Cells.Item(65, 9).Value = 3574.75   # I65: 0 -> 3574.75
$ws.Cells.Item(65, 10).Value = 5006.6   # J65: 111114980 -> 5006.6
$ws.Cells.Item(65, 11).Value = 17873.75   # K65: 0 -> 17873.75
$ws.Cells.Item(65, 12).Value = 25033   # L65: 555574900 -> 25033
$ws.Cells.Item(65, 13).Value = -14753.75   # M65: None -> -14753.75
$ws.Cells.Item(65, 14).Value = -31273   # N65: -555581140 -> -31273
# Row 81
$ws.Cells.Item(81, 8).Value = 47621020   # H81: 71430744 -> 47621020
$ws.Cells.Item(81, 9).Value = 71430410   # I81: 76925030 -> 71430410
$ws.Cells.Item(81, 10).Value = 2242.8572   # J81: 5000 -> 2242.8572
$ws.Cells.Item(81, 11).Value = 142860820   # K81: 153850060 -> 142860820
$ws.Cells.Item(81, 12).Value = 4485.7144   # L81: 10000 -> 4485.7144
$ws.Cells.Item(81, 13).Value = -142859759   # M81: -153848999 -> -142859759
$ws.Cells.Item(81, 14).Value = -6607.7144   # N81: -12122 -> -6607.7144
# Row 84
$ws.Cells.Item(84, 8).Value = 47621020   # H84: 71430744 -> 47621020
$ws.Cells.Item(84, 9).Value = 71430410   # I84: 76925030 -> 71430410
$ws.Cells.Item(84, 10).Value = 2242.8572   # J84: 5000 -> 2242.8572
$ws.Cells.Item(84, 11).Value = 714304100   # K84: 769250300 -> 714304100
$ws.Cells.Item(84, 12).Value = 22428.572   # L84: 50000 -> 22428.572
$ws.Cells.Item(84, 13).Value = -714298796   # M84: -769244996 -> -714298796
$ws.Cells.Item(84, 14).Value = -33036.572   # N84: -60608 -> -33036.572
# Row 123
$ws.Cells.Item(123, 8).Value = 44000   # H123: 46429 -> 44000
$ws.Cells.Item(123, 10).Value = 44000   # J123: 46429 -> 44000
$ws.Cells.Item(123, 12).Value = 44000   # L123: 46429 -> 44000
$ws.Cells.Item(123, 14).Value = -53800   # N123: -56229 -> -53800
# Row 136
$ws.Cells.Item(136, 8).Value = 1418.3928   # H136: 1423.9286 -> 1418.3928
$ws.Cells.Item(136, 9).Value = 1087.2174   # I136: 1122.091 -> 1087.2174
$ws.Cells.Item(136, 10).Value = 2941.8   # J136: 2530.6667 -> 2941.8
$ws.Cells.Item(136, 11).Value = 3261.6522   # K136: 3366.273 -> 3261.6522
$ws.Cells.Item(136, 12).Value = 8825.400000000001   # L136: 7592.000100000001 -> 8825.400000000001
$ws.Cells.Item(136, 13).Value = -711.6522   # M136: -816.2729999999997 -> -711.6522
$ws.Cells.Item(136, 14).Value = -13925.4   # N136: -12692.0001 -> -13925.4
